$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1. Bold the "Q.x" question paragraphs (adds <w:b/><w:bCs/> to the runs
#    and the paragraph mark).
# ----------------------------------------------------------------------
$boldIndices = @(4,6,8,10,12,14,16,18,20,22,24,26,28,30,32,34,36,38,40,42)
foreach ($i in $boldIndices) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    $r.Font.Bold = $true
    $r.Font.BoldBi = $true
}

# ----------------------------------------------------------------------
# 2. Append a second sentence as its own run onto the end of the
#    "Not directly, std::list doesn't support ..." answer (Q.7 answer).
#    Built by typing the continuation into a fresh paragraph then
#    folding that paragraph back into the previous one (deleting the
#    paragraph mark) so the two sentences stay separate runs without
#    leaving stray run-formatting behind.
# ----------------------------------------------------------------------
$q7AnswerIndex = 29
$p = $d.Paragraphs.Item($q7AnswerIndex)
$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($q7AnswerIndex + 1)
$p2.Range.InsertAfter(" This  because lists are a 2 linked dynamic array that don't have set locations for each element, so there is no list[1] to get the second element because there is no location 1, you have to check the first element and find its next link through that element.")
$p1b = $d.Paragraphs.Item($q7AnswerIndex)
$markEnd = $p1b.Range.End
$markRange = $d.Range($markEnd - 1, $markEnd)
$markRange.Delete()

# ----------------------------------------------------------------------
# 3. Fill in the empty paragraph under Q.10 with the answer, written as
#    two runs the same way.
# ----------------------------------------------------------------------
$q10AnswerIndex = 35
$p = $d.Paragraphs.Item($q10AnswerIndex)
$p.Range.InsertAfter("ParticleClass ")
$p2 = $d.Paragraphs.Item($q10AnswerIndex)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($q10AnswerIndex + 1)
$p3.Range.InsertAfter("Default constructor was called.")
$p2b = $d.Paragraphs.Item($q10AnswerIndex)
$markEnd2 = $p2b.Range.End
$markRange2 = $d.Range($markEnd2 - 1, $markEnd2)
$markRange2.Delete()
